# Update "想去人数" (number of people interested) values on the 展览 and
# 全部类型 sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 203
$ws1.Range("F4").Value = 806
$ws1.Range("F5").Value = 71
$ws1.Range("F6").Value = 22

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 203
$ws4.Range("F5").Value = 806
$ws4.Range("F6").Value = 71
$ws4.Range("F7").Value = 22
